# This script rotates the species-observation data among rows 2, 3 and 4:
#   old row 3 -> row 2
#   old row 4 -> row 3
#   old row 2 -> row 4
# Only the columns that actually vary between the three records change
# (A, B, E, F, G, H, Q, R, Z, AB); every other column already holds the
# same value in all three rows, so nothing else needs to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

# Capture the current ("before") values for the three rows.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value()
    $row3[$col] = $ws.Range($col + "3").Value()
    $row4[$col] = $ws.Range($col + "4").Value()
}

# Write the rotated values back: row3 -> row2, row4 -> row3, row2 -> row4.
foreach ($col in $cols) {
    $ws.Range($col + "2").Value = $row3[$col]
    $ws.Range($col + "3").Value = $row4[$col]
    $ws.Range($col + "4").Value = $row2[$col]
}
